$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.546.47'
$ws.Range('E2').Value = '  -0.64%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.442.24'
$ws.Range('E3').Value = '  -1.22%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '408.01'
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.52'
$ws.Range('E6').Value = '  +1.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.594'
$ws.Range('E7').Value = '  -1.69%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -1.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.126'
$ws.Range('E10').Value = '  -2.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.26'
$ws.Range('E11').Value = '  -2.28%  '
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.651.35'
$ws.Range('E13').Value = '  +5.16%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.47'
$ws.Range('E14').Value = '  -3.26%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '19.99'
$ws.Range('E15').Value = '  -1.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '62.534.73'
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '11.44'
$ws.Range('E17').Value = '  +5.52%  '
$ws.Range('E18').Value = '  -2.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000136'
$ws.Range('E19').Value = '  -3.75%  '
$ws.Range('E20').Value = '  -5.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '84.05'
$ws.Range('E21').Value = '  +1.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '314.29'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('E23').Value = '  -1.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.16'
$ws.Range('E24').Value = '  -0.74%  '
$ws.Range('E25').Value = '  +7.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.20'
$ws.Range('E27').Value = '  -1.25%  '
$ws.Range('E28').Value = '  +6.52%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.66'
$ws.Range('E29').Value = '  -1.26%  '
$ws.Range('E30').Value = '  -3.52%  '
$ws.Range('E31').Value = '  -3.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '42.48'
$ws.Range('E32').Value = '  -1.21%  '
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('E34').Value = '  -4.19%  '
$ws.Range('E35').Value = '  -2.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '51.60'
$ws.Range('E36').Value = '  -1.63%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.998'
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('E38').Value = '  -5.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.324'
$ws.Range('E39').Value = '  +12.52%  '
$ws.Range('E40').Value = '  -2.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '138.43'
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('E42').Value = '  -0.35%  '
$ws.Range('E43').Value = '  -0.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.00'
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.84'
$ws.Range('E45').Value = '  -4.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.22'
$ws.Range('E46').Value = '  -1.37%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '21.37'
$ws.Range('E47').Value = '  -4.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.129.17'
$ws.Range('E49').Value = '  -2.97%  '
$ws.Range('E50').Value = '  +2.93%  '
$ws.Range('E51').Value = '  +21.85%  '
